$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the obsolete BOM lines per the new 3.3V-only design:
#   - U2 (AAT3221IGV-3.3-T1 LDO regulator)            -> row 4
#   - Q1, Q2, Q3 (BSS138 N-channel MOSFETs)            -> row 9 (before prior deletes)
#   - R1-R6 (10K resistors)                            -> row 10 (before prior deletes)
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(4).Delete()

# The defined name's range shrinks by the one row removed from inside it.
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Sheet1!`$A`$1:`$E`$6"
}

# Reflect where the user's selection ended up after editing the table.
$ws.Range("B8").Select()
